$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("K1").Value = "N26"
$ws.Range("L1").Value = "N26sigma"

$ws.Range("K2").Value = 30000
$ws.Range("L2").Value = 200

$ws.Range("K3").Select()
